$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = "191,63 "; $ws.Cells.Item(2, 3).Value = "195,73 "; $ws.Cells.Item(2, 4).Value = "138,80 "
$ws.Cells.Item(3, 2).Value = "146,95 "; $ws.Cells.Item(3, 3).Value = "149,94 "; $ws.Cells.Item(3, 4).Value = "81,84 "
$ws.Cells.Item(4, 2).Value = "156,06 "; $ws.Cells.Item(4, 3).Value = "157,26 "; $ws.Cells.Item(4, 4).Value = "77,73 "
$ws.Cells.Item(5, 2).Value = "216,98 "; $ws.Cells.Item(5, 3).Value = "223,12 "; $ws.Cells.Item(5, 4).Value = "116,66 "
$ws.Cells.Item(6, 2).Value = "176,33 "; $ws.Cells.Item(6, 3).Value = "170,95 "; $ws.Cells.Item(6, 4).Value = "83,35 "
$ws.Cells.Item(7, 2).Value = "208,15 "; $ws.Cells.Item(7, 3).Value = "194,32 "; $ws.Cells.Item(7, 4).Value = "127,44 "
$ws.Cells.Item(8, 2).Value = "191,86 "; $ws.Cells.Item(8, 3).Value = "180,62 "; $ws.Cells.Item(8, 4).Value = "126,04 "
$ws.Cells.Item(9, 2).Value = "188,00 "; $ws.Cells.Item(9, 3).Value = "180,27 "; $ws.Cells.Item(9, 4).Value = "115,38 "
$ws.Cells.Item(10, 2).Value = "212,13 "; $ws.Cells.Item(10, 3).Value = "204,05 "; $ws.Cells.Item(10, 4).Value = "121,58 "
$ws.Cells.Item(11, 2).Value = "214,23 "; $ws.Cells.Item(11, 3).Value = "236,84 "; $ws.Cells.Item(11, 4).Value = "129,94 "
$ws.Cells.Item(12, 2).Value = "213,42 "; $ws.Cells.Item(12, 3).Value = "185,89 "; $ws.Cells.Item(12, 4).Value = "115,14 "
$ws.Cells.Item(13, 2).Value = "178,04 "; $ws.Cells.Item(13, 3).Value = "169,84 "; $ws.Cells.Item(13, 4).Value = "97,22 "
$ws.Cells.Item(14, 2).Value = "182,78 "; $ws.Cells.Item(14, 3).Value = "188,97 "; $ws.Cells.Item(14, 4).Value = "107,22 "
$ws.Cells.Item(15, 2).Value = "210,72 "; $ws.Cells.Item(15, 3).Value = "208,26 "; $ws.Cells.Item(15, 4).Value = "117,13 "
$ws.Cells.Item(16, 2).Value = "221,38 "; $ws.Cells.Item(16, 3).Value = "210,37 "; $ws.Cells.Item(16, 4).Value = "86,21 "
$ws.Cells.Item(17, 2).Value = "176,52 "; $ws.Cells.Item(17, 3).Value = "189,40 "; $ws.Cells.Item(17, 4).Value = "94,17 "
$ws.Cells.Item(18, 2).Value = "189,17 "; $ws.Cells.Item(18, 3).Value = "187,18 "; $ws.Cells.Item(18, 4).Value = "109,17 "
$ws.Cells.Item(19, 2).Value = "147,46 "; $ws.Cells.Item(19, 3).Value = "137,00 "; $ws.Cells.Item(19, 4).Value = "87,53 "
$ws.Cells.Item(20, 2).Value = "244,34 "; $ws.Cells.Item(20, 3).Value = "219,39 "; $ws.Cells.Item(20, 4).Value = "105,54 "
$ws.Cells.Item(21, 2).Value = "173,88 "; $ws.Cells.Item(21, 3).Value = "177,25 "; $ws.Cells.Item(21, 4).Value = "92,83 "
$ws.Cells.Item(22, 2).Value = "206,27 "; $ws.Cells.Item(22, 3).Value = "204,40 "; $ws.Cells.Item(22, 4).Value = "112,45 "
$ws.Cells.Item(23, 2).Value = "174,88 "; $ws.Cells.Item(23, 3).Value = "180,52 "; $ws.Cells.Item(23, 4).Value = "84,97 "
$ws.Cells.Item(24, 2).Value = "186,55 "; $ws.Cells.Item(24, 3).Value = "202,26 "; $ws.Cells.Item(24, 4).Value = "109,01 "
$ws.Cells.Item(25, 2).Value = "186,36 "; $ws.Cells.Item(25, 3).Value = "189,64 "; $ws.Cells.Item(25, 4).Value = "108,00 "
$ws.Cells.Item(26, 2).Value = "188,47 "; $ws.Cells.Item(26, 3).Value = "185,77 "; $ws.Cells.Item(26, 4).Value = "92,42 "
$ws.Cells.Item(27, 2).Value = "131,75 "; $ws.Cells.Item(27, 3).Value = "141,84 "; $ws.Cells.Item(27, 4).Value = "81,31 "
